$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.751.19'
$ws.Range("E2").Value = '  -2.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.715.12'
$ws.Range("E3").Value = '  -6.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.34'
$ws.Range("E5").Value = '  -4.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.77'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  -3.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.726.36'
$ws.Range("E9").Value = '  -6.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.10'
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -2.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.348'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.189.01'
$ws.Range("E14").Value = '  -6.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.634.62'
$ws.Range("E15").Value = '  -2.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.59'
$ws.Range("E16").Value = '  -4.53%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000135'
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.711.55'
$ws.Range("E18").Value = '  -6.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.72'
$ws.Range("E19").Value = '  -3.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.90'
$ws.Range("E20").Value = '  -5.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.36'
$ws.Range("E21").Value = '  -4.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.19'
$ws.Range("E22").Value = '  -5.97%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.60'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.75'
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.423'
$ws.Range("E26").Value = '  -5.48%  '
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0832'
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.42'
$ws.Range("E30").Value = '  -3.82%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  -3.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.03'
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.63'
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.18'
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.34'
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.937'
$ws.Range("E37").Value = '  -5.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.13'
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("E39").Value = '  -4.54%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.38'
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.175.12'
$ws.Range("E41").Value = '  -6.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.51'
$ws.Range("E42").Value = '  -3.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.995'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0551'
$ws.Range("E44").Value = '  -2.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.599'
$ws.Range("E45").Value = '  -6.71%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.78'
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.80'
$ws.Range("E47").Value = '  -9.06%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.35'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0226'
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0882'
$ws.Range("E50").Value = '  -4.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.89'
$ws.Range("E51").Value = '  -1.70%  '
